$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daylight Savings Time update: drop the fixed clock-times from the
# recurring reset labels (times now shift with DST so they're no longer
# accurate to hard-code).
$ws.Range("A2").Value = "Most Dailies"
$ws.Range("A3").Value = "Hilts Trader Limited Time Items"
$ws.Range("A4").Value = "Hilts Trader Limited Time Items"
$ws.Range("A5").Value = "Weekly Reset"

# Move the active selection to A5.
[void]$ws.Range("A5").Select()
